# Refresh the cryptos price table (GitHub Actions data-update commit).
# Only columns D (Price) and E (Volume 1h) change; every other cell is left alone.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price values parse as plain numbers (e.g. "1.008", "2.110").
# The source sheet stores Price as literal text (leading/trailing zeros such as
# "2.110" or "0.01920" are significant), so force those specific cells to Text
# format before writing, which stops Excel from silently re-typing them as numbers
# and dropping the trailing zeros. Cells whose new value is not a valid number
# (e.g. "27.098.90") are left on General format since Excel keeps those as text anyway.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D18", "D19", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row-by-row updates: Price (D) and Volume(1h) (E) per the refreshed feed.
$ws.Range("D2").Value = "27.098.90"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "1.823.92"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").Value = "311.52"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "0.4628"
$ws.Range("D8").Value = "0.3639"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "0.07296"
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("D10").Value = "0.8698"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "1.867.76"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "0.07618"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "92.54"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "0.000008650"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "27.367.66"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "5.198"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").Value = "2.092.21"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "151.81"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "1.862"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "18.25"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "2.110"
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("D29").Value = "116.24"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "5.085"
$ws.Range("E30").Value = "  -3.90%  "
$ws.Range("D31").Value = "0.08888"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "0.7354"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").Value = "4.447"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "2.573"
$ws.Range("E37").Value = "  +6.95%  "
$ws.Range("D38").Value = "0.05263"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "1.069"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("D40").Value = "0.01920"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").Value = "2.934"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("D42").Value = "7.174"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").Value = "0.5214"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").Value = "0.1633"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").Value = "8.265"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "103.97"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "10.13"
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").Value = "0.06250"
$ws.Range("E51").Value = "  -1.30%  "
